# Refactored naming of stores to sinks
#
# Adds a new "Cache" cylinder (can) shape to slide 2, mirroring the style of
# the other "can"-shaped (Cylinder) flowchart nodes already on the slide, and
# places it as the very first shape in the slide's shape tree (i.e. sent to
# the back of the z-order), matching the target OOXML diff:
#
#   <p:sp>
#     <p:nvSpPr>
#       <p:cNvPr id="19" name="Cylinder 18"/>
#       ...
#     <p:spPr>
#       <a:xfrm><a:off x="4010158" y="5051350"/><a:ext cx="570449" cy="379570"/></a:xfrm>
#       <a:prstGeom prst="can"><a:avLst/></a:prstGeom>
#     ...
#     <p:txBody> ... "Cache" ...

# PowerPoint's Shape.Left/Top/Width/Height are COM `Single` (32-bit float)
# properties, so a plain EMU/12700 conversion can be off by 1 EMU once it
# round-trips through 32-bit float storage (because PowerPoint truncates when
# converting back to EMU on save). This helper nudges the point value by the
# smallest amount needed so that, after being narrowed to `single` precision,
# it converts back to exactly the target EMU value.
function Get-EmuSafePt($targetEmu) {
    $pt = $targetEmu / 12700.0
    $result = $pt
    for ($i = 0; $i -lt 5000; $i++) {
        $candidate = $pt + ($i * 0.0000001)
        $asSingle = [single]$candidate
        $backEmu = [math]::Floor([double]$asSingle * 12700.0)
        if ($backEmu -eq $targetEmu) {
            $result = $candidate
            break
        }
    }
    $result
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate an existing "can" (Cylinder) shape to use as a style template, so
# the new shape inherits the identical <p:style> block (lnRef/fillRef/
# effectRef/fontRef) and plain (no extra fill/line overrides) <p:spPr> used
# by the diagram's other cylinder shapes.
$template = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Cylinder 5") {
        $template = $s.Shapes.Item($i)
        break
    }
}

$newShape = $template.Duplicate().Item(1)
$newShape.Name = "Cylinder 18"

$newShape.Left = Get-EmuSafePt 4010158
$newShape.Top = Get-EmuSafePt 5051350
$newShape.Width = Get-EmuSafePt 570449
$newShape.Height = Get-EmuSafePt 379570

$newShape.TextFrame.TextRange.Text = "Cache"

# Move the new shape to the very start of the shape tree (bottom of z-order),
# matching its position immediately after </p:grpSpPr> in the target diff.
$newShape.ZOrder(1)
